$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 9: "Subscribe Event" service spec row ---
$ws.Range("A9").Value = "Subscribe Event"
$ws.Range("B9").Value = "<application-context>/studentcenter/subscribe-event"
$ws.Range("C9").Value = "POST"
$ws.Range("D9").Value = "Content-Type: application/json"
$ws.Range("E9").Value = "{`n ""eventId"": ""14"",`n ""userId"": ""N01""`n}"
$ws.Range("F9").Value = "Both are mandatory"

# Borders around the whole new row (matches the bordered-table look used elsewhere)
$ws.Range("A9:F9").Borders.LineStyle = 1

# E9 (body) wraps like the other "Body" column cells
$ws.Range("E9").WrapText = $true

# Row height to match the rest of the sheet's multi-line rows
$ws.Rows.Item(9).RowHeight = 60

# Restore the selection to the cell below the new data, like the source workbook
[void]$ws.Range("B8").Select()
